$d = $word.ActiveDocument

# 1. Fix screenshot temp location: Program.sc.ImageFiles -> Model.ImageFiles
$d.Content.Find.Execute("Program.sc.ImageFiles", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Model.ImageFiles", 2)

# 2. Word re-seats the hidden "_GoBack" bookmark on the most recent edit: it now
#    brackets the point right after "Model" (inside the word we just typed).
$r = $d.Content
$r.Find.Execute("Model", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$goBack = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $goBack)

# 3. Shrink the page margins from 1" to 0.5" on every side.
$ps = $d.PageSetup
$ps.TopMargin = 36
$ps.BottomMargin = 36
$ps.LeftMargin = 36
$ps.RightMargin = 36
